$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "dataInvestigate_file_inspector"
$ws.Range("B25").Value = "File Inspector"
$ws.Range("A26").Value = "files"
$ws.Range("B26").Value = "Files"

$ws.Range("A26").Select()
